# Modify combo damage scaling (onblock value for "slk" move in the raw data
# table) and adjust the mirrored/derived table accordingly. Formulas in the
# sheet recalc automatically, cascading the change to D5/G5 (first table)
# and G25 (blockstun formula in the second table). Finally, move the active
# selection to G25 as recorded in the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Raw input data lives in the second table (rows 23-38). Row 25 corresponds
# to the "slk" move; column D is "onblock".
$ws.Range("D25").Value = -12

# Make sure all dependent formulas (D5, G5, G25, ...) are recalculated.
$excel.Calculate()

# Update the saved selection/active cell for the sheet view.
$ws.Activate()
$ws.Range("G25").Select()
